# "update table from cruise to month"
#
# A new summary sheet, "count_rda_statistics", is inserted in front of the
# existing "count_rda_axis" / "count_rda_margin" sheets and becomes the
# active tab. It reports the overall R-squared / Adj. R-squared of the RDA
# model (full model vs. the backward-selected model).

$wb = $excel.ActiveWorkbook

# Insert the new worksheet before the current first sheet so it becomes
# sheet #1 and the active tab (mirrors Excel's own Worksheets.Add(Before:=...)).
$firstSheet = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($firstSheet)
$ws.Name = "count_rda_statistics"

# Header row - bold + centered, matching the other report sheets.
$ws.Range("A1:C1").Font.Bold = $true
$ws.Range("A1:C1").HorizontalAlignment = -4108  # xlCenter

$ws.Range("A1").Value = "Model"
$ws.Range("B1").Value = "R.squared"
$ws.Range("C1").Value = "Adj.R.squared"

$ws.Range("A2").Value = "Full model"
$ws.Range("B2").Value = 0.3824413253767041
$ws.Range("C2").Value = 0.2095248964821812

$ws.Range("A3").Value = "Backward selected"
$ws.Range("B3").Value = 0.3143555125444172
$ws.Range("C3").Value = 0.1873843111637536

$ws.Activate()
